$wb = $excel.ActiveWorkbook

# --- Add the new "CatalogPricerule" worksheet as the last tab ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "CatalogPricerule"

# --- Header row (row 1) ---
$ws.Range("A1").Value = "DataSet"
$ws.Range("B1").Value = "UserName"
$ws.Range("C1").Value = "Password"
$ws.Range("D1").Value = "RuleName"
$ws.Range("E1").Value = "Description"
$ws.Range("F1").Value = "Financecategory"
$ws.Range("G1").Value = "Attribute"
$ws.Range("H1").Value = "HFCategory"
$ws.Range("I1").Value = "HFsubcategory"
$ws.Range("J1").Value = "OxoCategory"
$ws.Range("K1").Value = "Discount"
$ws.Range("L1").Value = "Oxosubcategory"
$ws.Range("M1").Value = "URL"
$ws.Range("N1").Value = "preprodURL"
$ws.Range("O1").Value = "pageTitle"
$ws.Range("P1").Value = "Priority"
$ws.Range("Q1").Value = "Discard subsequent rules"
$ws.Range("R1").Value = "Apply"
$ws.Range("A1:R1").Interior.Color = 65535

# --- Row 2 ---
$ws.Range("A2").Value = "AccountDetails"
$ws.Range("B2").Value = "mkoppanadam@helenoftroy.com"
$ws.Range("C2").Value = "Amtlmcflmipq1!"

# --- Row 3 ---
$ws.Range("A3").Value = "Catalogpricedetails"
$ws.Range("D3").Value = "Qatestcatalogrule"
$ws.Range("E3").Value = "5%offcatalogpricerule"
$ws.Range("F3").Value = "TRADE"
$ws.Range("G3").Value = "Category"
$ws.Range("P3").Value = 0

# --- Row 4 ---
$ws.Range("A4").Value = "Categoryselection"
$ws.Range("H4").Value = "Bottles & Drinkware"
$ws.Range("I4").Value = "Coffee"
$ws.Range("J4").Value = "Coffee & Beverage"
$ws.Range("K4").Value = 5
$ws.Range("L4").Value = "Coffee & Tea"
$ws.Range("Q4").Value = "Yes"
$ws.Range("R4").Value = "Apply as percentage of original"

# --- Row 5 ---
$ws.Range("A5").Value = "Hydroflask"
$ws.Range("M5").Value = "https://mcloud-na-stage.hydroflask.com/"
$ws.Range("N5").Value = "https://mcloud-na-preprod.hydroflask.com/"
$ws.Range("O5").Value = "Home Page (Hydroflask)"

# --- Row 6 ---
$ws.Range("A6").Value = "OXO"
$ws.Range("M6").Value = "https://mcloud-na-stage.oxo.com/"
$ws.Range("N6").Value = "https://mcloud-na-preprod.oxo.com/"
$ws.Range("O6").Value = "Home Page (OXO)"

# --- Column widths to (roughly) match the authored layout ---
$ws.Columns.Item(1).ColumnWidth = 18.28515625
$ws.Columns.Item(2).ColumnWidth = 31.5703125
$ws.Columns.Item(3).ColumnWidth = 15.5703125
$ws.Columns.Item(4).ColumnWidth = 16.85546875
$ws.Columns.Item(5).ColumnWidth = 20.7109375
$ws.Columns.Item(6).ColumnWidth = 15.5703125
$ws.Columns.Item(8).ColumnWidth = 19
$ws.Columns.Item(9).ColumnWidth = 14.140625
$ws.Columns.Item(10).ColumnWidth = 17.85546875
$ws.Columns.Item(12).ColumnWidth = 15.42578125
$ws.Columns.Item(13).ColumnWidth = 38.7109375
$ws.Columns.Item(14).ColumnWidth = 41.140625
$ws.Columns.Item(15).ColumnWidth = 22.7109375
$ws.Columns.Item(16).ColumnWidth = 7.5703125
$ws.Columns.Item(17).ColumnWidth = 23.5703125
$ws.Columns.Item(18).ColumnWidth = 28.85546875

# --- View: scroll the new sheet so column J is leftmost, select M9 ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 10
$ws.Range("M9").Select()
